# This script restructures the "114_1" nominations summary sheet:
# - removes the old top-level "Summary"/"New nominations" style rows
# - relabels each category's sub-rows with the category name baked in
#   (e.g. "New nominations" -> "Civilian, New nominations")
# - adds "Returned to White House" for Civilian, and two new summary rows
#   ("Total new nominations" / "Total carryover nominations") replacing
#   the old "Summary" header row, shrinking the sheet from 40 to 39 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the last row (old row 40, "Total returned to the White House") since
# the new layout only needs 39 rows; its content is re-emitted at row 39 below.
$ws.Rows(40).Delete()

# Cell B38 previously held a "#,##0"-formatted value (old row 38, 2207) but
# in the new layout it holds "Total withdrawn " (335) which must go back to
# the plain general-number format used elsewhere (same format as B39).
# Copy that formatting over before the values below overwrite both cells.
$ws.Range("B39").Copy()
$ws.Range("B38").PasteSpecial(-4122)

$ws.Range("A1").Value = 'Labels'
$ws.Range("B1").Value = 'Values'
$ws.Range("A2").Value = 'Congress'
$ws.Range("B2").Value = 114
$ws.Range("A3").Value = 'Session'
$ws.Range("B3").Value = 1
$ws.Range("A4").Value = 'Start Date'
$ws.Range("B4").Value = 42010
$ws.Range("A5").Value = 'End Date'
$ws.Range("B5").Value = 42369
$ws.Range("A6").Value = 'Civilian'
$ws.Range("A7").Value = '     Civilian, New nominations'
$ws.Range("B7").Value = 366
$ws.Range("A8").Value = '     Civilian, Confirmed '
$ws.Range("B8").Value = 173
$ws.Range("A9").Value = '     Civilian, Unconfirmed '
$ws.Range("B9").Value = 181
$ws.Range("A10").Value = '     Civilian, Withdrawn '
$ws.Range("B10").Value = 10
$ws.Range("A11").Value = '     Civilian, Returned to White House '
$ws.Range("B11").Value = 2
$ws.Range("A12").Value = 'Other Civilian'
$ws.Range("A13").Value = '     Other Civilian, New nominations'
$ws.Range("B13").Value = 3802
$ws.Range("B13").NumberFormat = "#,##0"
$ws.Range("A14").Value = '     Other Civilian, Confirmed '
$ws.Range("B14").Value = 3383
$ws.Range("B14").NumberFormat = "#,##0"
$ws.Range("A15").Value = '     Other Civilian, Unconfirmed '
$ws.Range("B15").Value = 97
$ws.Range("A16").Value = '     Other Civilian, Withdrawn '
$ws.Range("B16").Value = 322
$ws.Range("A17").Value = 'Air Force'
$ws.Range("A18").Value = '     Air Force, New nominations'
$ws.Range("B18").Value = 5734
$ws.Range("B18").NumberFormat = "#,##0"
$ws.Range("A19").Value = '     Air Force, Confirmed '
$ws.Range("B19").Value = 5550
$ws.Range("B19").NumberFormat = "#,##0"
$ws.Range("A20").Value = '     Air Force, Unconfirmed '
$ws.Range("B20").Value = 181
$ws.Range("A21").Value = '     Air Force, Withdrawn '
$ws.Range("B21").Value = 3
$ws.Range("A22").Value = 'Army'
$ws.Range("A23").Value = '     Army, New nominations'
$ws.Range("B23").Value = 5214
$ws.Range("B23").NumberFormat = "#,##0"
$ws.Range("A24").Value = '     Army, Confirmed '
$ws.Range("B24").Value = 3474
$ws.Range("B24").NumberFormat = "#,##0"
$ws.Range("A25").Value = '     Army, Unconfirmed '
$ws.Range("B25").Value = 1740
$ws.Range("B25").NumberFormat = "#,##0"
$ws.Range("A26").Value = 'Navy'
$ws.Range("A27").Value = '     Navy, New nominations'
$ws.Range("B27").Value = 3936
$ws.Range("B27").NumberFormat = "#,##0"
$ws.Range("A28").Value = '     Navy, Confirmed '
$ws.Range("B28").Value = 3931
$ws.Range("B28").NumberFormat = "#,##0"
$ws.Range("A29").Value = '     Navy, Unconfirmed '
$ws.Range("B29").Value = 5
$ws.Range("A30").Value = 'Marine Corps'
$ws.Range("A31").Value = '     Marine Corps, New nominations'
$ws.Range("B31").Value = 1070
$ws.Range("B31").NumberFormat = "#,##0"
$ws.Range("A32").Value = '     Marine Corps, Confirmed '
$ws.Range("B32").Value = 1067
$ws.Range("B32").NumberFormat = "#,##0"
$ws.Range("A33").Value = '     Marine Corps, Unconfirmed '
$ws.Range("B33").Value = 3
$ws.Range("A34").Value = 'Total new nominations'
$ws.Range("B34").Value = 20122
$ws.Range("B34").NumberFormat = "#,##0"
$ws.Range("A35").Value = 'Total carryover nominations'
$ws.Range("B35").Value = 0
$ws.Range("A36").Value = 'Total confirmed '
$ws.Range("B36").Value = 17578
$ws.Range("B36").NumberFormat = "#,##0"
$ws.Range("A37").Value = 'Total unconfirmed '
$ws.Range("B37").Value = 2207
$ws.Range("B37").NumberFormat = "#,##0"
$ws.Range("A38").Value = 'Total withdrawn '
$ws.Range("B38").Value = 335
$ws.Range("A39").Value = 'Total returned to the White House '
$ws.Range("B39").Value = 2
